$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Several paragraphs in the table had their sentence split across multiple
#    <w:r> runs (an interleaved bookmark, or <w:proofErr/> "grammar" markers,
#    used to sit in the middle of the sentence). Re-typing the identical full
#    sentence over the already-found text collapses it back down to a single
#    run and drops whatever markup used to sit in between.
# ---------------------------------------------------------------------------

$mergedSentences = @(
  "Böylelikle bugün izinli olan çalışanları görüp işlerimi onlara göre planlayabilirim.",
  "İzin isteyen çalışanın adını , soyadını, departmanını ,izin bakiyesini, izin talep nedenini, izin istediği aralığı",
  "Şirkette bulunan departmanları, çalışan sayısını ve hangi departmanda hangi çalışanların bulunduğunu görebilmeliyim.",
  "Böylelikle çalışanları ve departmanları daha kolay yönetebilir , çalışanların departmanlar arası geçişini kolaylıkla yapabilirim.",
  "Diğer Talepler.."
)

foreach ($sentence in $mergedSentences) {
    $rng = $d.Content
    $rng.Find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, $sentence, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) The document's last table cell gets a brand new, empty paragraph
#    appended after "Böylelikle planlamalarımı rahatlıkla yapabilirim.", and
#    the "_GoBack" bookmark (previously sitting mid-sentence in the first
#    merged paragraph above) now lives alone in that new empty paragraph.
# ---------------------------------------------------------------------------

$tbl = $d.Tables.Item(1)
$lastRow = $tbl.Rows.Count
$cell = $tbl.Cell($lastRow, 1)
$cellEnd = $cell.Range
$cellEnd.Collapse(0)   # wdCollapseEnd
$cellEnd.InsertParagraphAfter()

# Re-find the anchor sentence to read back a trustworthy absolute offset for
# the paragraph mark that now separates it from the freshly inserted, empty
# paragraph right after it.
$locator = $d.Content
$locator.Find.Execute("Böylelikle planlamalarımı rahatlıkla yapabilirim.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraMarkPos = $locator.End

$bookRange = $d.Range($paraMarkPos, $paraMarkPos + 2)
$d.Bookmarks.Add("_GoBack", $bookRange) | Out-Null
